$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI Lamc2-Itgb1 results per Dr Hou advice: added ECs as a sending/target cluster
# and refreshed all computed statistics for the FAPs/sCs/ECs x FAPs/sCs/ECs grid.

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Lamc2"
$ws.Range("C2").Value = "Itgb1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2121133333333333
$ws.Range("H2").Value = 0.6363399999999999
$ws.Range("I2").Value = 0.04296779043029777
$ws.Range("J2").Value = 0.04296779043029776
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 112.513392
$ws.Range("N2").Value = 337.540176
$ws.Range("O2").Value = 0.3275312977368564
$ws.Range("P2").Value = 0.3275312977368564
$ws.Range("Q2").Value = 23.86559062176
$ws.Range("R2").Value = 214.79031559584
$ws.Range("S2").Value = 0.01407329616052071
$ws.Range("T2").Value = 0.0140732961605207

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Lamc2"
$ws.Range("C3").Value = "Itgb1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2121133333333333
$ws.Range("H3").Value = 0.6363399999999999
$ws.Range("I3").Value = 0.04296779043029777
$ws.Range("J3").Value = 0.04296779043029776
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 106.314466
$ws.Range("N3").Value = 318.943398
$ws.Range("O3").Value = 0.3094859589441663
$ws.Range("P3").Value = 0.3094859589441664
$ws.Range("Q3").Value = 22.55071576481333
$ws.Range("R3").Value = 202.95644188332
$ws.Range("S3").Value = 0.01329792782503268
$ws.Range("T3").Value = 0.01329792782503268

# Row 4: ECs -> sCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Lamc2"
$ws.Range("C4").Value = "Itgb1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2121133333333333
$ws.Range("H4").Value = 0.6363399999999999
$ws.Range("I4").Value = 0.04296779043029777
$ws.Range("J4").Value = 0.04296779043029776
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 124.6916553333333
$ws.Range("N4").Value = 374.074966
$ws.Range("O4").Value = 0.3629827433189773
$ws.Range("P4").Value = 0.3629827433189773
$ws.Range("Q4").Value = 26.44876265160444
$ws.Range("R4").Value = 238.03886386444
$ws.Range("S4").Value = 0.01559656644474438
$ws.Range("T4").Value = 0.01559656644474438

# Row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Lamc2"
$ws.Range("C5").Value = "Itgb1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.865410333333333
$ws.Range("H5").Value = 11.596231
$ws.Range("I5").Value = 0.7830160344930734
$ws.Range("J5").Value = 0.7830160344930733
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 112.513392
$ws.Range("N5").Value = 337.540176
$ws.Range("O5").Value = 0.3275312977368564
$ws.Range("P5").Value = 0.3275312977368564
$ws.Range("Q5").Value = 434.910428075184
$ws.Range("R5").Value = 3914.193852676656
$ws.Range("S5").Value = 0.2564622579262834
$ws.Range("T5").Value = 0.2564622579262834

# Row 6: FAPs -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Lamc2"
$ws.Range("C6").Value = "Itgb1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.865410333333333
$ws.Range("H6").Value = 11.596231
$ws.Range("I6").Value = 0.7830160344930734
$ws.Range("J6").Value = 0.7830160344930733
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 106.314466
$ws.Range("N6").Value = 318.943398
$ws.Range("O6").Value = 0.3094859589441663
$ws.Range("P6").Value = 0.3094859589441664
$ws.Range("Q6").Value = 410.9490354592153
$ws.Range("R6").Value = 3698.541319132938
$ws.Range("S6").Value = 0.2423324683037472
$ws.Range("T6").Value = 0.2423324683037472

# Row 7: FAPs -> sCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Lamc2"
$ws.Range("C7").Value = "Itgb1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.865410333333333
$ws.Range("H7").Value = 11.596231
$ws.Range("I7").Value = 0.7830160344930734
$ws.Range("J7").Value = 0.7830160344930733
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 124.6916553333333
$ws.Range("N7").Value = 374.074966
$ws.Range("O7").Value = 0.3629827433189773
$ws.Range("P7").Value = 0.3629827433189773
$ws.Range("Q7").Value = 481.9844130059051
$ws.Range("R7").Value = 4337.859717053146
$ws.Range("S7").Value = 0.2842213082630428
$ws.Range("T7").Value = 0.2842213082630428

# Row 8: sCs -> ECs
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Lamc2"
$ws.Range("C8").Value = "Itgb1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.8590423333333334
$ws.Range("H8").Value = 2.577127
$ws.Range("I8").Value = 0.1740161750766289
$ws.Range("J8").Value = 0.1740161750766288
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 112.513392
$ws.Range("N8").Value = 337.540176
$ws.Range("O8").Value = 0.3275312977368564
$ws.Range("P8").Value = 0.3275312977368564
$ws.Range("Q8").Value = 96.653766794928
$ws.Range("R8").Value = 869.8839011543519
$ws.Range("S8").Value = 0.05699574365005225
$ws.Range("T8").Value = 0.05699574365005224

# Row 9: sCs -> FAPs
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Lamc2"
$ws.Range("C9").Value = "Itgb1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.8590423333333334
$ws.Range("H9").Value = 2.577127
$ws.Range("I9").Value = 0.1740161750766289
$ws.Range("J9").Value = 0.1740161750766288
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 106.314466
$ws.Range("N9").Value = 318.943398
$ws.Range("O9").Value = 0.3094859589441663
$ws.Range("P9").Value = 0.3094859589441664
$ws.Range("Q9").Value = 91.32862693972733
$ws.Range("R9").Value = 821.957642457546
$ws.Range("S9").Value = 0.05385556281538641
$ws.Range("T9").Value = 0.05385556281538642

# Row 10: sCs -> sCs
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Lamc2"
$ws.Range("C10").Value = "Itgb1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.8590423333333334
$ws.Range("H10").Value = 2.577127
$ws.Range("I10").Value = 0.1740161750766289
$ws.Range("J10").Value = 0.1740161750766288
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 124.6916553333333
$ws.Range("N10").Value = 374.074966
$ws.Range("O10").Value = 0.3629827433189773
$ws.Range("P10").Value = 0.3629827433189773
$ws.Range("Q10").Value = 107.1154105447425
$ws.Range("R10").Value = 964.038694902682
$ws.Range("S10").Value = 0.0631648686111902
$ws.Range("T10").Value = 0.06316486861119018
